# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets hold the same set of rows (2-19) for this section of data,
# and the same F-column counts were updated on both.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8882
    3  = 8343
    5  = 166
    6  = 209
    7  = 254
    8  = 758
    9  = 220
    10 = 5485
    11 = 9
    14 = 22
    17 = 165
    18 = 216
    19 = 17
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
